$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-16 Saturday" "2024-11-17 Sunday"

Replace-Text "23×82=" "55×47="
Replace-Text "27×78=" "37×11="
Replace-Text "13×52=" "86×49="
Replace-Text "18×33=" "81×38="
Replace-Text "16×74=" "70×20="

Replace-Text "59×13=" "54×12="
Replace-Text "85×25=" "36×41="
Replace-Text "95×92=" "11×52="
Replace-Text "77×66=" "89×82="
Replace-Text "80×99=" "60×68="

Replace-Text "70×40=" "13×19="
Replace-Text "49×39=" "87×86="
Replace-Text "59×68=" "48×99="
Replace-Text "95×11=" "44×68="
Replace-Text "72×88=" "69×36="

Replace-Text "83×68=" "78×71="
Replace-Text "11×75=" "71×29="
Replace-Text "89×59=" "23×50="
Replace-Text "75×30=" "32×26="
Replace-Text "36×29=" "58×79="

Replace-Text "40×68=" "49×43="
Replace-Text "74×85=" "60×48="
Replace-Text "74×48=" "27×78="
Replace-Text "12×11=" "77×52="
Replace-Text "79×71=" "22×23="

Write-Output "done"
